# Actualización automática 2025-08-05 10:15:08
# Update budget ("PRESUPUESTO") figures on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 1500
$ws.Range("G4").Value = 300
$ws.Range("G5").Value = 1000
$ws.Range("G6").Value = 8000
$ws.Range("G8").Value = 1150
$ws.Range("G11").Value = 1500
$ws.Range("G12").Value = 6000
$ws.Range("G13").Value = 2000
$ws.Range("G14").Value = 0
$ws.Range("G16").Value = 1200
$ws.Range("G19").Value = 1000
$ws.Range("G20").Value = 2000

# G24 is a total row; recompute explicitly to match the expected literal value
$ws.Range("G24").Value = 48450
